$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.649.37'
$ws.Range('E2').Value = '  -3.72%  '
$ws.Range('D3').Value = '2.511.43'
$ws.Range('E3').Value = '  -5.28%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'575.70"
$ws.Range('E5').Value = '  -2.71%  '
$ws.Range('D6').Value = "'167.12"
$ws.Range('E6').Value = '  -4.77%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').Value = '2.509.95'
$ws.Range('E9').Value = '  -5.28%  '
$ws.Range('E10').Value = '  -7.24%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = "'0.343"
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '2.970.85'
$ws.Range('E14').Value = '  -5.24%  '
$ws.Range('D15').Value = '69.489.00'
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('E16').Value = '  -6.30%  '
$ws.Range('D17').Value = "'24.91"
$ws.Range('E17').Value = '  -4.46%  '
$ws.Range('D18').Value = '2.513.32'
$ws.Range('E18').Value = '  -5.63%  '
$ws.Range('D19').Value = "'11.44"
$ws.Range('E19').Value = '  -7.01%  '
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').Value = "'350.63"
$ws.Range('E21').Value = '  -5.70%  '
$ws.Range('E22').Value = '  -5.12%  '
$ws.Range('D23').Value = "'1.97"
$ws.Range('E23').Value = '  -5.47%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = "'68.81"
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').Value = "'4.02"
$ws.Range('E26').Value = '  -6.34%  '
$ws.Range('E27').Value = '  -7.48%  '
$ws.Range('D28').Value = '2.640.60'
$ws.Range('E28').Value = '  -5.04%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  -6.40%  '
$ws.Range('D31').Value = "'7.91"
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('D32').Value = "'477.22"
$ws.Range('E32').Value = '  -5.33%  '
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').Value = '  -3.51%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = "'154.49"
$ws.Range('E37').Value = '  -4.47%  '
$ws.Range('D38').Value = "'18.92"
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('D39').Value = "'18.55"
$ws.Range('E39').Value = '  -4.30%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').Value = "'4.75"
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('E42').Value = '  -2.98%  '
$ws.Range('E43').Value = '  -7.82%  '
$ws.Range('D44').Value = "'1.18"
$ws.Range('E44').Value = '  -13.24%  '
$ws.Range('E45').Value = '  -9.87%  '
$ws.Range('D46').Value = "'38.13"
$ws.Range('E46').Value = '  -2.57%  '
$ws.Range('D47').Value = "'144.36"
$ws.Range('E47').Value = '  -6.18%  '
$ws.Range('D48').Value = "'0.531"
$ws.Range('E48').Value = '  -4.00%  '
$ws.Range('E49').Value = '  -3.85%  '
$ws.Range('E50').Value = '  -5.43%  '
$ws.Range('D51').Value = "'0.0731"
$ws.Range('E51').Value = '  -2.72%  '
